$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared string for row 18's label
$ws.Range("A18").Value = "main cuts with 8.5 ns upper delta t cut"
$ws.Range("B18").Value = -17.692699999999999
$ws.Range("C18").Formula = "=B18-B14"
$ws.Range("D18").Value = -17.821300000000001
$ws.Range("E18").Formula = "=D18-D14"
$ws.Range("F18").Value = 2.5070000000000001
$ws.Range("G18").Value = 0.17299999999999999

$ws.Range("A19").Value = "main cuts with E = 500 MeV cut"
$ws.Range("B19").Value = -17.692799999999998
$ws.Range("C19").Formula = "=B19-B14"
$ws.Range("D19").Value = -17.821400000000001
$ws.Range("E19").Formula = "=D19-D14"
$ws.Range("F19").Value = 2.2530000000000001
$ws.Range("G19").Value = 0.155

$ws.Range("D30").Select()
